# Updates cryptos list values (price/volume columns) per the commit diff.
# Some "Price" (column D) values are plain decimal numbers (e.g. "577.45") that
# Excel would silently convert to a numeric cell type if assigned directly,
# dropping significant trailing zeros / rewriting scientific notation.
# For those cells we assign with a leading apostrophe (forces text entry) and
# then reset the cell style back to "Normal" so no stray number-format/style
# is left behind, matching the original unstyled inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.132.56"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "3.584.62"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'577.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Value = "'188.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.09%  "
$ws.Range("E7").Value = "  -2.45%  "
$ws.Range("D8").Value = "3.580.39"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "'0.660"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  -3.51%  "
$ws.Range("D13").Value = "'0.0000301"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "'9.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.19%  "
$ws.Range("D15").Value = "4.151.25"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("E16").Value = "  +2.97%  "
$ws.Range("D17").Value = "3.576.38"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "69.972.82"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  -1.17%  "
$ws.Range("D22").Value = "'474.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.59%  "
$ws.Range("D23").Value = "'19.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.65%  "
$ws.Range("D24").Value = "'5.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -8.59%  "
$ws.Range("D25").Value = "'4.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "'88.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").Value = "'11.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").Value = "'32.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "
$ws.Range("D31").Value = "'7.67"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").Value = "'12.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").Value = "'65.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "'577.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.55%  "
$ws.Range("D36").Value = "'38.86"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "0.0₃0799"
$ws.Range("E38").Value = "  -4.19%  "
$ws.Range("D39").Value = "'0.396"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E40").Value = "  -5.31%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'3.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.42%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'3.21"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.18%  "
$ws.Range("D43").Value = "'2.90"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.43%  "
$ws.Range("D44").Value = "3.233.66"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").Value = "'3.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.40%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "'9.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.10%  "
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  -3.74%  "
